$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: header spacer row. Copy formatting from S2 into new T2 cell ---
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)

# --- Row 3: year headers. Q3:S3 change style to match D3:P3 (s=11), add T3 = 2023 ---
$ws.Range("P3").Copy()
$ws.Range("Q3:T3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# --- Row 4: data values. Copy S4 formatting into T4, set new value ---
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 263951

# --- Row 5: data values. Copy S5 formatting into T5, set new value ---
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = 3.7

# --- Row 6: data values. Copy S6 formatting into T6, set new value ---
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = 32.299999999999997

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 16.5
$ws.Rows.Item(3).RowHeight = 15.75

# --- Remove the lingering cell selection marker left in sheetView ---
$ws.Range("A1").Select()

$excel.CutCopyMode = $false
